$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, report week dates) ---
$ws.Range("A8").Value = "Volume 30   Number  50"
$ws.Range("C9").Value = "Report Covering the Week  12/11/2023  Through  12/17/2023"

# --- Cells changing from numeric to placeholder text (shared strings for 0 / ***.*), reuse format from a same-style neighbor ---
$ws.Range("D15").Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4122) | Out-Null
$ws.Range("C15").Formula = "=TEXT(0,""0"")"
$ws.Range("C15").Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4163) | Out-Null

$ws.Range("D22").Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4122) | Out-Null
$ws.Range("C22").Formula = "=TEXT(0,""0"")"
$ws.Range("C22").Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4163) | Out-Null

$ws.Range("D15").Copy() | Out-Null
$ws.Range("D26").PasteSpecial(-4122) | Out-Null
$ws.Range("D26").Formula = "=TEXT(0,""0"")"
$ws.Range("D26").Copy() | Out-Null
$ws.Range("D26").PasteSpecial(-4163) | Out-Null

$ws.Range("E15").Copy() | Out-Null
$ws.Range("E26").PasteSpecial(-4122) | Out-Null
$ws.Range("E26").Formula = "=TEXT(0,""***.*"")"
$ws.Range("E26").Copy() | Out-Null
$ws.Range("E26").PasteSpecial(-4163) | Out-Null

$ws.Range("G30").Copy() | Out-Null
$ws.Range("F30").PasteSpecial(-4122) | Out-Null
$ws.Range("F30").Formula = "=TEXT(0,""0"")"
$ws.Range("F30").Copy() | Out-Null
$ws.Range("F30").PasteSpecial(-4163) | Out-Null

# --- Cells changing from placeholder text to numeric value, reuse format from a same-style/col neighbor ---
$ws.Range("C20").Copy() | Out-Null
$ws.Range("D20").PasteSpecial(-4122) | Out-Null
$ws.Range("D20").Value = 4
$ws.Range("K16").Copy() | Out-Null
$ws.Range("E20").PasteSpecial(-4122) | Out-Null
$ws.Range("E20").Value = 50

# --- Simple numeric value updates ---
$ws.Range("M15").Value = 30.769230769230
$ws.Range("N15").Value = -50
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 14
$ws.Range("H16").Value = 16.666666666666
$ws.Range("I16").Value = 173
$ws.Range("J16").Value = 197
$ws.Range("K16").Value = -12.182741116751
$ws.Range("L16").Value = 50.434782608695
$ws.Range("M16").Value = -33.969465648855
$ws.Range("N16").Value = -85.314091680814
$ws.Range("C17").Value = 6
$ws.Range("E17").Value = -40
$ws.Range("F17").Value = 22
$ws.Range("H17").Value = -24.137931034482
$ws.Range("I17").Value = 360
$ws.Range("J17").Value = 370
$ws.Range("K17").Value = -2.702702702702
$ws.Range("L17").Value = 39.534883720930
$ws.Range("M17").Value = 54.506437768240
$ws.Range("N17").Value = 0.840336134453
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 13
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 151
$ws.Range("K18").Value = -11.920529801324
$ws.Range("L18").Value = 23.148148148148
$ws.Range("M18").Value = -62.108262108262
$ws.Range("N18").Value = -90.627202255109
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = -25
$ws.Range("G19").Value = 32
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 331
$ws.Range("J19").Value = 370
$ws.Range("K19").Value = -10.540540540540
$ws.Range("L19").Value = 16.140350877193
$ws.Range("M19").Value = -3.498542274052
$ws.Range("N19").Value = -40.787119856887
$ws.Range("F20").Value = 21
$ws.Range("G20").Value = 18
$ws.Range("H20").Value = 16.666666666666
$ws.Range("I20").Value = 303
$ws.Range("J20").Value = 250
$ws.Range("K20").Value = 21.2
$ws.Range("L20").Value = 33.480176211453
$ws.Range("M20").Value = 8.214285714285
$ws.Range("N20").Value = -90.992865636147
$ws.Range("D21").Value = 30
$ws.Range("E21").Value = -20
$ws.Range("F21").Value = 106
$ws.Range("H21").Value = 0.952380952380
$ws.Range("I21").Value = 1324
$ws.Range("J21").Value = 1363
$ws.Range("K21").Value = -2.861335289801
$ws.Range("L21").Value = 29.931305201177
$ws.Range("M21").Value = -10.961667787491
$ws.Range("N21").Value = -80.872580179139
$ws.Range("L22").Value = 55.555555555555
$ws.Range("C24").Value = 21
$ws.Range("D24").Value = 52
$ws.Range("E24").Value = -59.615384615384
$ws.Range("F24").Value = 96
$ws.Range("G24").Value = 161
$ws.Range("H24").Value = -40.372670807453
$ws.Range("I24").Value = 1324
$ws.Range("J24").Value = 1615
$ws.Range("K24").Value = -18.018575851393
$ws.Range("L24").Value = 30.830039525691
$ws.Range("M24").Value = 71.725032425421
$ws.Range("C25").Value = 13
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = 44.444444444444
$ws.Range("F25").Value = 59
$ws.Range("G25").Value = 38
$ws.Range("H25").Value = 55.263157894736
$ws.Range("I25").Value = 556
$ws.Range("J25").Value = 525
$ws.Range("K25").Value = 5.904761904761
$ws.Range("L25").Value = 11.2
$ws.Range("M25").Value = -7.794361525704
$ws.Range("F26").Value = 8
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = 166.666666666667
$ws.Range("I26").Value = 37
$ws.Range("K26").Value = 2.777777777777
$ws.Range("L26").Value = 12.121212121212
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = -50
$ws.Range("F27").Value = 3
$ws.Range("H27").Value = -40
$ws.Range("I27").Value = 52
$ws.Range("J27").Value = 52
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = -24.637681159420
$ws.Range("N28").Value = -60.606060606060
$ws.Range("N29").Value = -67.857142857142

Write-Host "done"
